$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (existing rows 59..118 shift down to 60..119,
# carrying their original content and formatting with them).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record.
$ws.Range("A59").Value = 5
$ws.Range("B59").Value = "Macroferia Regional de Talca"
$ws.Range("C59").Value = "Maule"
$ws.Range("D59").Value = 45167
$ws.Range("E59").Value = 7
$ws.Range("F59").Value = 100112040
$ws.Range("G59").Value = "Cilantro"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 150
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 10000
$ws.Range("M59").Value = 10000
$ws.Range("N59").Value = '$/caja 36 atados'
$ws.Range("O59").Value = "Región Metropolitana"
$ws.Range("P59").Value = 278
$ws.Range("Q59").Value = 36
$ws.Range("R59").Value = "Hortaliza"

# Make sure the new row's date cell carries the same date style as the rest
# of column D (the Insert() above should already propagate it from row 58/60,
# but set it explicitly to be safe).
$ws.Range("D59").NumberFormat = $ws.Range("D60").NumberFormat
